$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")
Write-Host $ws.Name
